$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume text values (which often look numeric, e.g. "0.9990")
# are written back as literal text, matching the original inline-string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.138.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.835.45'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.17'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6654'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2956'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07357'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.78'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07681'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.840.03'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.018'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6752'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.27'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.197'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.145.90'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008246'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.91'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.20%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9993'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.304'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '161.03'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1422'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.683'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.01'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.502'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.230'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.105'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.204'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05322'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.19%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7461'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.130'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.679'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.316.94'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01805'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.715'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9255'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.996'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9987'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.52'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.985.49'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5169'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.276'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07436'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +9.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05928'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.01%  '
